# added decoupling caps to neopixles and Arduino
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row so the old "47 uF Ceramic Cap" row (row 9) moves down to row 10,
# leaving room for the new decoupling cap line item at row 9.
$ws.Rows("9:9").Insert()

# Row 8 used to be the generic "10 uF Cap" (P975-ND). Replace it with the new
# 100 uF electrolytic decoupling cap used on the board.
$ws.Range("B8").Value() = "100 uF Electrolytic Cap"
$ws.Range("C8").Value() = "732-8598-1-ND"
$ws.Range("F8").Value() = 0.1

# New row 9: 1uF ceramic decoupling caps added for the neopixels and Arduino.
$ws.Range("B9").Value() = "1uF Ceramic Cap"
$ws.Range("C9").Value() = "1276-1182-1-ND"
$ws.Range("D9").Value() = 3
$ws.Range("E9").Value() = 5
$ws.Range("F9").Value() = 0.1
$ws.Range("G9").Formula() = "=E9*F9"

# Row 10 (previously row 9, "47 uF Ceramic Cap"): quantities updated and the
# "NEW" note cleared now that the part isn't new anymore.
$ws.Range("D10").Value() = 1
$ws.Range("E10").Value() = 2
$ws.Range("H10").ClearContents()

# Leave the selection where the author last left it.
[void]$ws.Range("F9").Select()
